# "Generate Report for handoff" for 9c0760a9-f347-488d-8db0-4bdee644a390.md
#
# The report generation marks that file's latest handoff as freshly
# kicked off: its status flips from "Handed back" to "Not yet handed off"
# (surfaced on the Overview sheet and on each per-language detail sheet),
# and the per-language "Latest Handoff Datetime" is stamped with the
# handoff time just recorded.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 9c0760a9-...md file; B = zh-cn status, C = de-de status
$overview.Range("B3").Value = "Not yet handed off"
$overview.Range("C3").Value = "Not yet handed off"

# zh-cn detail sheet: row 3 is the same file; B = Status, D = Latest Handoff Datetime
$zhcn.Range("B3").Value = "Not yet handed off"
$zhcn.Range("D3").Value = "2016-01-07 07:47:14"

# de-de detail sheet: row 3 is the same file; B = Status, D = Latest Handoff Datetime
$dede.Range("B3").Value = "Not yet handed off"
$dede.Range("D3").Value = "2016-01-07 07:47:24"
